$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the new columns.
#    Original layout (A..P):
#      A modalidade, B autoria_classificacao, C total, D total_sucesso,
#      E particip, F taxa_sucesso, G arrecadado_sucesso, H media_sucesso,
#      I std_sucesso, J min_sucesso, K max_sucesso, L apoio_medio,
#      M contribuicoes, N media_contribuicoes, O menor_ano, P maior_ano
#
#    Insert 3 columns right after L (apoio_medio) for apoio_std/min/max,
#    then (after that shift) insert 3 more columns right after the column
#    that now holds media_contribuicoes (contribuicoes_med) for the new
#    contribuicoes std/min/max columns.
# ---------------------------------------------------------------------------
$ws.Columns.Item(13).Insert()
$ws.Columns.Item(13).Insert()
$ws.Columns.Item(13).Insert()

$ws.Columns.Item(18).Insert()
$ws.Columns.Item(18).Insert()
$ws.Columns.Item(18).Insert()

# ---------------------------------------------------------------------------
# 2. Rename headers to the new column names.
#    After the inserts the layout is:
#      A modalidade, B autoria_classificacao, C total, D total_sucesso,
#      E particip, F taxa_sucesso, G arrecadado_sucesso, H media_sucesso,
#      I std_sucesso, J min_sucesso, K max_sucesso, L apoio_medio,
#      M (new), N (new), O (new),
#      P contribuicoes, Q media_contribuicoes,
#      R (new), S (new), T (new),
#      U menor_ano, V maior_ano
# ---------------------------------------------------------------------------
$ws.Range("H1").Value = "arrecadado_avg"
$ws.Range("I1").Value = "arrecadado_std"
$ws.Range("J1").Value = "arrecadado_min"
$ws.Range("K1").Value = "arrecadado_max"

$ws.Range("M1").Value = "apoio_std"
$ws.Range("N1").Value = "apoio_min"
$ws.Range("O1").Value = "apoio_max"

$ws.Range("Q1").Value = "contribuicoes_med"
$ws.Range("R1").Value = "contribuicoes_std"
$ws.Range("S1").Value = "contribuicoes_min"
$ws.Range("T1").Value = "contribuicoes_max"

# ---------------------------------------------------------------------------
# 3. Update data rows 2-6.
#    L (apoio_medio) gets a recomputed value; M/N/O (apoio_std/min/max) are
#    brand-new columns; P/Q (contribuicoes / contribuicoes_med) keep their
#    previous values, and R/S/T (contribuicoes std/min/max) are brand new.
# ---------------------------------------------------------------------------

# Row 2 - coletivo
$ws.Range("L2").Value = 31.16847126718795
$ws.Range("M2").Value = 26.9469146898807
$ws.Range("N2").Value = 8.140546434454963
$ws.Range("O2").Value = 84.0771316599004
$ws.Range("R2").Value = 2.497617912751115
$ws.Range("S2").Value = 3
$ws.Range("T2").Value = 10

# Row 3 - empresa
$ws.Range("L3").Value = 42.14013096402113
$ws.Range("M3").Value = 8.830628986869351
$ws.Range("N3").Value = 35.89593332526331
$ws.Range("O3").Value = 48.38432860277894
$ws.Range("R3").Value = 3.535533905932738
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = 15

# Row 4 - feminino
$ws.Range("L4").Value = 19.42257389357928
$ws.Range("M4").Value = 8.876706034650423
$ws.Range("N4").Value = 5.929916345397809
$ws.Range("O4").Value = 35.80030877323957
$ws.Range("R4").Value = 20.56688435388656
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 79

# Row 5 - masculino
$ws.Range("L5").Value = 17.81312171425239
$ws.Range("M5").Value = 9.871079671113662
$ws.Range("N5").Value = 6.098311514417047
$ws.Range("O5").Value = 45.46067338136409
$ws.Range("R5").Value = 45.19114957599552
$ws.Range("S5").Value = 1
$ws.Range("T5").Value = 208

# Row 6 - outros
$ws.Range("L6").Value = 21.37695663886886
$ws.Range("M6").Value = 15.58070588764584
$ws.Range("N6").Value = 1.011042153300025
$ws.Range("O6").Value = 70.01644246718027
$ws.Range("R6").Value = 30.79468667274807
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = 196
